$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. These cells hold numeric-looking text
# (shared strings) in the original workbook, so we force a Text number
# format before assignment to stop Excel from auto-converting the
# string into a float, then clear the format override back to the
# sheet's default (General) so no stray per-cell formatting remains.
$updates = @{
    "B13" = "69.58"
    "D13" = "80.69"
    "B14" = "24.18"
    "C14" = "44.87"
    "D14" = "69.05"
    "B16" = "85.95"
    "C16" = "13.71"
    "D16" = "99.66"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
